# "Fruta, Feria Lagunitas de Puerto Montt - Damasco" - weekly refresh
# Inserts 3 new daily entries (2022-01-11, variety "Modesto") at the top of the
# data block and 2 more (2022-01-07) lower down, shifting the remaining historical
# rows down by 3 (old row 33 becomes row 36). Net effect: dimension grows from
# A1:T33 to A1:T36.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Cells.Item(21, 1).Value = 4
$ws.Cells.Item(21, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(21, 3).Value = "Los Lagos"
$ws.Cells.Item(21, 4).Value = 44572
$ws.Cells.Item(21, 5).Value = 10
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100103
$ws.Cells.Item(21, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(21, 9).Value = 100103003
$ws.Cells.Item(21, 10).Value = "Damasco"
$ws.Cells.Item(21, 11).Value = "Modesto"
$ws.Cells.Item(21, 12).Value = "Especial"
$ws.Cells.Item(21, 13).Value = 150
$ws.Cells.Item(21, 14).Value = 21000
$ws.Cells.Item(21, 15).Value = 21000
$ws.Cells.Item(21, 16).Value = 21000
$ws.Cells.Item(21, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(21, 18).Value = "Región Metropolitana"
$ws.Cells.Item(21, 19).Value = 1167
$ws.Cells.Item(21, 20).Value = 18

# Row 22
$ws.Cells.Item(22, 1).Value = 4
$ws.Cells.Item(22, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(22, 3).Value = "Los Lagos"
$ws.Cells.Item(22, 4).Value = 44572
$ws.Cells.Item(22, 5).Value = 10
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100103
$ws.Cells.Item(22, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(22, 9).Value = 100103003
$ws.Cells.Item(22, 10).Value = "Damasco"
$ws.Cells.Item(22, 11).Value = "Modesto"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 150
$ws.Cells.Item(22, 14).Value = 18000
$ws.Cells.Item(22, 15).Value = 18000
$ws.Cells.Item(22, 16).Value = 18000
$ws.Cells.Item(22, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(22, 18).Value = "Región Metropolitana"
$ws.Cells.Item(22, 19).Value = 1000
$ws.Cells.Item(22, 20).Value = 18

# Row 23
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 44572
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100103
$ws.Cells.Item(23, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(23, 9).Value = 100103003
$ws.Cells.Item(23, 10).Value = "Damasco"
$ws.Cells.Item(23, 11).Value = "Modesto"
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 150
$ws.Cells.Item(23, 14).Value = 16000
$ws.Cells.Item(23, 15).Value = 16000
$ws.Cells.Item(23, 16).Value = 16000
$ws.Cells.Item(23, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(23, 18).Value = "Región Metropolitana"
$ws.Cells.Item(23, 19).Value = 889
$ws.Cells.Item(23, 20).Value = 18

# Row 24
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(24, 3).Value = "Los Lagos"
$ws.Cells.Item(24, 4).Value = 44553
$ws.Cells.Item(24, 5).Value = 10
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100103
$ws.Cells.Item(24, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(24, 9).Value = 100103003
$ws.Cells.Item(24, 10).Value = "Damasco"
$ws.Cells.Item(24, 11).Value = "Castle Brite"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 250
$ws.Cells.Item(24, 14).Value = 20000
$ws.Cells.Item(24, 15).Value = 20000
$ws.Cells.Item(24, 16).Value = 20000
$ws.Cells.Item(24, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(24, 18).Value = "Región Metropolitana"
$ws.Cells.Item(24, 19).Value = 1111
$ws.Cells.Item(24, 20).Value = 18

# Row 25
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(25, 3).Value = "Los Lagos"
$ws.Cells.Item(25, 4).Value = 44553
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100103
$ws.Cells.Item(25, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(25, 9).Value = 100103003
$ws.Cells.Item(25, 10).Value = "Damasco"
$ws.Cells.Item(25, 11).Value = "Castle Brite"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 250
$ws.Cells.Item(25, 14).Value = 18000
$ws.Cells.Item(25, 15).Value = 18000
$ws.Cells.Item(25, 16).Value = 18000
$ws.Cells.Item(25, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(25, 18).Value = "Región Metropolitana"
$ws.Cells.Item(25, 19).Value = 1000
$ws.Cells.Item(25, 20).Value = 18

# Row 26
$ws.Cells.Item(26, 1).Value = 4
$ws.Cells.Item(26, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(26, 3).Value = "Los Lagos"
$ws.Cells.Item(26, 4).Value = 44553
$ws.Cells.Item(26, 5).Value = 10
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100103
$ws.Cells.Item(26, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(26, 9).Value = 100103003
$ws.Cells.Item(26, 10).Value = "Damasco"
$ws.Cells.Item(26, 11).Value = "Castle Brite"
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 250
$ws.Cells.Item(26, 14).Value = 16000
$ws.Cells.Item(26, 15).Value = 16000
$ws.Cells.Item(26, 16).Value = 16000
$ws.Cells.Item(26, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(26, 18).Value = "Región Metropolitana"
$ws.Cells.Item(26, 19).Value = 889
$ws.Cells.Item(26, 20).Value = 18

# Row 27
$ws.Cells.Item(27, 1).Value = 4
$ws.Cells.Item(27, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(27, 3).Value = "Los Lagos"
$ws.Cells.Item(27, 4).Value = 44551
$ws.Cells.Item(27, 5).Value = 10
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100103
$ws.Cells.Item(27, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(27, 9).Value = 100103003
$ws.Cells.Item(27, 10).Value = "Damasco"
$ws.Cells.Item(27, 11).Value = "Castle Brite"
$ws.Cells.Item(27, 12).Value = "Especial"
$ws.Cells.Item(27, 13).Value = 200
$ws.Cells.Item(27, 14).Value = 20000
$ws.Cells.Item(27, 15).Value = 20000
$ws.Cells.Item(27, 16).Value = 20000
$ws.Cells.Item(27, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(27, 18).Value = "Región Metropolitana"
$ws.Cells.Item(27, 19).Value = 1111
$ws.Cells.Item(27, 20).Value = 18

# Row 28
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(28, 3).Value = "Los Lagos"
$ws.Cells.Item(28, 4).Value = 44551
$ws.Cells.Item(28, 5).Value = 10
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100103
$ws.Cells.Item(28, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(28, 9).Value = 100103003
$ws.Cells.Item(28, 10).Value = "Damasco"
$ws.Cells.Item(28, 11).Value = "Castle Brite"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 14).Value = 18000
$ws.Cells.Item(28, 15).Value = 18000
$ws.Cells.Item(28, 16).Value = 18000
$ws.Cells.Item(28, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(28, 18).Value = "Región Metropolitana"
$ws.Cells.Item(28, 19).Value = 1000
$ws.Cells.Item(28, 20).Value = 18

# Row 29
$ws.Cells.Item(29, 1).Value = 4
$ws.Cells.Item(29, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(29, 3).Value = "Los Lagos"
$ws.Cells.Item(29, 4).Value = 44551
$ws.Cells.Item(29, 5).Value = 10
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100103
$ws.Cells.Item(29, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(29, 9).Value = 100103003
$ws.Cells.Item(29, 10).Value = "Damasco"
$ws.Cells.Item(29, 11).Value = "Castle Brite"
$ws.Cells.Item(29, 12).Value = "Segunda"
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 14).Value = 16000
$ws.Cells.Item(29, 15).Value = 16000
$ws.Cells.Item(29, 16).Value = 16000
$ws.Cells.Item(29, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(29, 18).Value = "Región Metropolitana"
$ws.Cells.Item(29, 19).Value = 889
$ws.Cells.Item(29, 20).Value = 18

# Row 30
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value = "Los Lagos"
$ws.Cells.Item(30, 4).Value = 44187
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100103
$ws.Cells.Item(30, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(30, 9).Value = 100103003
$ws.Cells.Item(30, 10).Value = "Damasco"
$ws.Cells.Item(30, 11).Value = "Castle Brite"
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 350
$ws.Cells.Item(30, 14).Value = 16000
$ws.Cells.Item(30, 15).Value = 16000
$ws.Cells.Item(30, 16).Value = 16000
$ws.Cells.Item(30, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(30, 18).Value = "Región Metropolitana"
$ws.Cells.Item(30, 19).Value = 1067
$ws.Cells.Item(30, 20).Value = 15

# Row 31
$ws.Cells.Item(31, 1).Value = 4
$ws.Cells.Item(31, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(31, 3).Value = "Los Lagos"
$ws.Cells.Item(31, 4).Value = 44187
$ws.Cells.Item(31, 5).Value = 10
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100103
$ws.Cells.Item(31, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(31, 9).Value = 100103003
$ws.Cells.Item(31, 10).Value = "Damasco"
$ws.Cells.Item(31, 11).Value = "Castle Brite"
$ws.Cells.Item(31, 12).Value = "Segunda"
$ws.Cells.Item(31, 13).Value = 300
$ws.Cells.Item(31, 14).Value = 13000
$ws.Cells.Item(31, 15).Value = 13000
$ws.Cells.Item(31, 16).Value = 13000
$ws.Cells.Item(31, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(31, 18).Value = "Región Metropolitana"
$ws.Cells.Item(31, 19).Value = 867
$ws.Cells.Item(31, 20).Value = 15

# Row 32
$ws.Cells.Item(32, 1).Value = 4
$ws.Cells.Item(32, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(32, 3).Value = "Los Lagos"
$ws.Cells.Item(32, 4).Value = 44568
$ws.Cells.Item(32, 5).Value = 10
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100103
$ws.Cells.Item(32, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(32, 9).Value = 100103003
$ws.Cells.Item(32, 10).Value = "Damasco"
$ws.Cells.Item(32, 11).Value = "Castle Brite"
$ws.Cells.Item(32, 12).Value = "Especial"
$ws.Cells.Item(32, 13).Value = 200
$ws.Cells.Item(32, 14).Value = 21000
$ws.Cells.Item(32, 15).Value = 21000
$ws.Cells.Item(32, 16).Value = 21000
$ws.Cells.Item(32, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(32, 18).Value = "Región Metropolitana"
$ws.Cells.Item(32, 19).Value = 1167
$ws.Cells.Item(32, 20).Value = 18

# Row 33
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44568
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100103
$ws.Cells.Item(33, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(33, 9).Value = 100103003
$ws.Cells.Item(33, 10).Value = "Damasco"
$ws.Cells.Item(33, 11).Value = "Castle Brite"
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 200
$ws.Cells.Item(33, 14).Value = 18000
$ws.Cells.Item(33, 15).Value = 18000
$ws.Cells.Item(33, 16).Value = 18000
$ws.Cells.Item(33, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(33, 18).Value = "Región Metropolitana"
$ws.Cells.Item(33, 19).Value = 1000
$ws.Cells.Item(33, 20).Value = 18

# Row 34
$ws.Cells.Item(34, 1).Value = 4
$ws.Cells.Item(34, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(34, 3).Value = "Los Lagos"
$ws.Cells.Item(34, 4).Value = 44568
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 10
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100103
$ws.Cells.Item(34, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(34, 9).Value = 100103003
$ws.Cells.Item(34, 10).Value = "Damasco"
$ws.Cells.Item(34, 11).Value = "Castle Brite"
$ws.Cells.Item(34, 12).Value = "Segunda"
$ws.Cells.Item(34, 13).Value = 200
$ws.Cells.Item(34, 14).Value = 16000
$ws.Cells.Item(34, 15).Value = 16000
$ws.Cells.Item(34, 16).Value = 16000
$ws.Cells.Item(34, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(34, 18).Value = "Región Metropolitana"
$ws.Cells.Item(34, 19).Value = 889
$ws.Cells.Item(34, 20).Value = 18

# Row 35
$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(35, 3).Value = "Los Lagos"
$ws.Cells.Item(35, 4).Value = 44194
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 10
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100103
$ws.Cells.Item(35, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(35, 9).Value = 100103003
$ws.Cells.Item(35, 10).Value = "Damasco"
$ws.Cells.Item(35, 11).Value = "Castle Brite"
$ws.Cells.Item(35, 12).Value = "Segunda"
$ws.Cells.Item(35, 13).Value = 300
$ws.Cells.Item(35, 14).Value = 15000
$ws.Cells.Item(35, 15).Value = 16000
$ws.Cells.Item(35, 16).Value = 15500
$ws.Cells.Item(35, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(35, 18).Value = "Región Metropolitana"
$ws.Cells.Item(35, 19).Value = 1033
$ws.Cells.Item(35, 20).Value = 15

# Row 36
$ws.Cells.Item(36, 1).Value = 4
$ws.Cells.Item(36, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(36, 3).Value = "Los Lagos"
$ws.Cells.Item(36, 4).Value = 44540
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 10
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100103
$ws.Cells.Item(36, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(36, 9).Value = 100103003
$ws.Cells.Item(36, 10).Value = "Damasco"
$ws.Cells.Item(36, 11).Value = "Castle Brite"
$ws.Cells.Item(36, 12).Value = "Segunda"
$ws.Cells.Item(36, 13).Value = 600
$ws.Cells.Item(36, 14).Value = 16000
$ws.Cells.Item(36, 15).Value = 16000
$ws.Cells.Item(36, 16).Value = 16000
$ws.Cells.Item(36, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(36, 18).Value = "Región del Maule"
$ws.Cells.Item(36, 19).Value = 889
$ws.Cells.Item(36, 20).Value = 18

